$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.091.05'
$ws.Range("D3").Value = '2.928.28'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.11'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.20'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.506'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.89'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.61%  '
$ws.Range("E10").Value = '  +0.81%  '
$ws.Range("E12").Value = '  +1.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.78'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("E14").Value = '  -0.41%  '
$ws.Range("D15").Value = '3.411.57'
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("D16").Value = '61.038.34'
$ws.Range("E16").Value = '  +0.84%  '
$ws.Range("E17").Value = '  -1.48%  '
$ws.Range("D18").Value = '2.928.15'
$ws.Range("E18").Value = '  +0.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '432.33'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.82%  '
$ws.Range("E20").Value = '  -1.16%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '81.38'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.09'
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.22'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.01'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.44%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("E28").Value = '  +6.52%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.10'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.51'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("E33").Value = '  +1.05%  '
$ws.Range("D34").Value = '0.0₃0855'
$ws.Range("E34").Value = '  +1.99%  '
$ws.Range("E35").Value = '  +0.42%  '
$ws.Range("E36").Value = '  -0.39%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.06'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +3.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.125'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.18%  '
$ws.Range("E39").Value = '  -1.61%  '
$ws.Range("E40").Value = '  -1.55%  '
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.80'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '375.94'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.70%  '
$ws.Range("D44").Value = '2.728.15'
$ws.Range("E44").Value = '  +2.79%  '
$ws.Range("E45").Value = '  -0.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '130.67'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.32%  '
$ws.Range("E48").Value = '  -3.27%  '
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("E51").Value = '  +3.11%  '
